$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.028333286578941
$ws.Range("D2").Value = 1.036506160765282
$ws.Range("E2").Value = 1.028308089439664
$ws.Range("F2").Value = 1.045316223443975
$ws.Range("I2").Value = 1.033753740601323
$ws.Range("J2").Value = 1.033486186819577
$ws.Range("K2").Value = 1.039299707786021
$ws.Range("L2").Value = 1.031125284585217
$ws.Range("M2").Value = 1.048084805741589
$ws.Range("N2").Value = 1.005712725503983

$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.029264745061004
$ws.Range("D3").Value = 1.037196341351082
$ws.Range("E3").Value = 1.029098810820422
$ws.Range("F3").Value = 1.046146191492865
$ws.Range("I3").Value = 1.033896067206445
$ws.Range("J3").Value = 1.034058302321531
$ws.Range("K3").Value = 1.039799820788631
$ws.Range("L3").Value = 1.031723960519062
$ws.Range("M3").Value = 1.048726141269683

$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.029868095796984
$ws.Range("D4").Value = 1.037643300497118
$ws.Range("E4").Value = 1.029611384564347
$ws.Range("F4").Value = 1.046683933459939
$ws.Range("I4").Value = 1.033986983326679
$ws.Range("J4").Value = 1.034428524607848
$ws.Range("K4").Value = 1.040123116547517
$ws.Range("L4").Value = 1.032111613505976
$ws.Range("M4").Value = 1.049141176486108

$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.030121895208278
$ws.Range("D5").Value = 1.037831288166892
$ws.Range("E5").Value = 1.029827090249437
$ws.Range("F5").Value = 1.046910165357188
$ws.Range("I5").Value = 1.034024921664882
$ws.Range("J5").Value = 1.034584170967199
$ws.Range("K5").Value = 1.040258954475734
$ws.Range("L5").Value = 1.032274646155722
$ws.Range("M5").Value = 1.049315667406149

$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.030164518051775
$ws.Range("D6").Value = 1.037862857122597
$ws.Range("E6").Value = 1.029863321020171
$ws.Range("F6").Value = 1.046948160324017
$ws.Range("I6").Value = 1.034031275077548
$ws.Range("J6").Value = 1.034610304930425
$ws.Range("K6").Value = 1.040281757792864
$ws.Range("L6").Value = 1.032302023733621
$ws.Range("M6").Value = 1.049344965750858

$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.029871486486259
$ws.Range("D7").Value = 1.037645812061496
$ws.Range("E7").Value = 1.029614265973019
$ws.Range("F7").Value = 1.046686955734765
$ws.Range("I7").Value = 1.033987491372582
$ws.Range("J7").Value = 1.034430604344264
$ws.Range("K7").Value = 1.040124931919165
$ws.Range("L7").Value = 1.032113791708887
$ws.Range("M7").Value = 1.049143508004379

$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.028647945160686
$ws.Range("D8").Value = 1.03673933374085
$ws.Range("E8").Value = 1.028575125147773
$ws.Range("F8").Value = 1.04559656964865
$ws.Range("I8").Value = 1.033802084234763
$ws.Range("J8").Value = 1.033679530003825
$ws.Range("K8").Value = 1.039468786945269
$ws.Range("L8").Value = 1.03132755336186
$ws.Range("M8").Value = 1.048301537233928

$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.026496824634731
$ws.Range("D9").Value = 1.035144878708707
$ws.Range("E9").Value = 1.02675117283253
$ws.Range("F9").Value = 1.043680585440731
$ws.Range("I9").Value = 1.033466372774382
$ws.Range("J9").Value = 1.03235628549819
$ws.Range("K9").Value = 1.038310248527629
$ws.Range("L9").Value = 1.029944219626492
$ws.Range("M9").Value = 1.046818310057032

$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.025066120818482
$ws.Range("D10").Value = 1.034083937938518
$ws.Range("E10").Value = 1.025540102070435
$ws.Range("F10").Value = 1.042407005970271
$ws.Range("I10").Value = 1.033236550216357
$ws.Range("J10").Value = 1.031474354797171
$ws.Range("K10").Value = 1.037536396374287
$ws.Range("L10").Value = 1.029023493683477
$ws.Range("M10").Value = 1.045829868174028

$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.024447425944714
$ws.Range("D11").Value = 1.033625040336938
$ws.Range("E11").Value = 1.025016876093686
$ws.Range("F11").Value = 1.041856441358891
$ws.Range("I11").Value = 1.033135615480893
$ws.Range("J11").Value = 1.031092538374198
$ws.Range("K11").Value = 1.037200970650405
$ws.Range("L11").Value = 1.028625179244212
$ws.Range("M11").Value = 1.045401969031023

$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.024217737869335
$ws.Range("D12").Value = 1.033454661584798
$ws.Range("E12").Value = 1.024822704860062
$ws.Range("F12").Value = 1.041652074644083
$ws.Range("I12").Value = 1.033097911072045
$ws.Range("J12").Value = 1.030950725705319
$ws.Range("K12").Value = 1.037076328313773
$ws.Range("L12").Value = 1.028477283713681
$ws.Range("M12").Value = 1.045243045094109

$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.024267001167957
$ws.Range("D13").Value = 1.033491204929085
$ws.Range("E13").Value = 1.024864347160223
$ws.Range("F13").Value = 1.041695905748492
$ws.Range("I13").Value = 1.03310600841377
$ws.Range("J13").Value = 1.030981144490531
$ws.Range("K13").Value = 1.037103066774633
$ws.Range("L13").Value = 1.028509005245097
$ws.Range("M13").Value = 1.045277134032358

$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.024428437346015
$ws.Range("D14").Value = 1.033610955211317
$ws.Range("E14").Value = 1.025000822194873
$ws.Range("F14").Value = 1.041839545527647
$ws.Range("I14").Value = 1.033132503159379
$ws.Range("J14").Value = 1.031080815864831
$ws.Range("K14").Value = 1.037190668694957
$ws.Range("L14").Value = 1.028612953000354
$ws.Range("M14").Value = 1.045388831976699

$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.024527919803405
$ws.Range("D15").Value = 1.033684747464628
$ws.Range("E15").Value = 1.025084932626687
$ws.Range("F15").Value = 1.041928065007397
$ws.Range("I15").Value = 1.033148799272371
$ws.Range("J15").Value = 1.031142228172944
$ws.Range("K15").Value = 1.037244636501228
$ws.Range("L15").Value = 1.028677006124453
$ws.Range("M15").Value = 1.045457655043056

$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.025107198709705
$ws.Range("D16").Value = 1.034114404065151
$ws.Range("E16").Value = 1.025574851758338
$ws.Range("F16").Value = 1.042443564317492
$ws.Range("I16").Value = 1.03324321905404
$ws.Range("J16").Value = 1.031499696169543
$ws.Range("K16").Value = 1.037558650357454
$ws.Range("L16").Value = 1.029049936314349
$ws.Range("M16").Value = 1.045858268716639

$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.025470782556309
$ws.Range("D17").Value = 1.034384050436876
$ws.Range("E17").Value = 1.025882480945102
$ws.Range("F17").Value = 1.042767166580031
$ws.Range("I17").Value = 1.033302066318483
$ws.Range("J17").Value = 1.031723944701691
$ws.Range("K17").Value = 1.037755532059627
$ws.Range("L17").Value = 1.029283964400411
$ws.Range("M17").Value = 1.046109591595864

$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.025682932872437
$ws.Range("D18").Value = 1.034541378424488
$ws.Range("E18").Value = 1.026062029174009
$ws.Range("F18").Value = 1.042956005383517
$ws.Range("I18").Value = 1.033336253878044
$ws.Range("J18").Value = 1.031854751310113
$ws.Range("K18").Value = 1.037870336686479
$ws.Range("L18").Value = 1.029420504253573
$ws.Range("M18").Value = 1.046256193794191

$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.025755283844468
$ws.Range("D19").Value = 1.034595031264578
$ws.Range("E19").Value = 1.026123269655155
$ws.Range("F19").Value = 1.043020409253058
$ws.Range("I19").Value = 1.033347887687591
$ws.Range("J19").Value = 1.0318993540063
$ws.Range("K19").Value = 1.03790947644029
$ws.Range("L19").Value = 1.029467066753982
$ws.Range("M19").Value = 1.046306182976705

$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.0254317653732
$ws.Range("D20").Value = 1.034355114975053
$ws.Range("E20").Value = 1.025849463511308
$ws.Range("F20").Value = 1.042732438100388
$ws.Range("I20").Value = 1.033295766735809
$ws.Range("J20").Value = 1.031699884303541
$ws.Range("K20").Value = 1.037734411935386
$ws.Range("L20").Value = 1.029258851748887
$ws.Range("M20").Value = 1.046082626003916

$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.024380895017502
$ws.Range("D21").Value = 1.033575689620054
$ws.Range("E21").Value = 1.024960628767242
$ws.Range("F21").Value = 1.041797243383413
$ws.Range("I21").Value = 1.033124706985448
$ws.Range("J21").Value = 1.03105146481799
$ws.Range("K21").Value = 1.037164873491488
$ws.Range("L21").Value = 1.028582341416548
$ws.Range("M21").Value = 1.045355939239309

$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.023720881875455
$ws.Range("D22").Value = 1.033086075919961
$ws.Range("E22").Value = 1.024402814815714
$ws.Range("F22").Value = 1.041210045708527
$ws.Range("I22").Value = 1.033015924296492
$ws.Range("J22").Value = 1.030643841669302
$ws.Range("K22").Value = 1.03680649170505
$ws.Range("L22").Value = 1.028157318198017
$ws.Range("M22").Value = 1.044899140663326

$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.024070699570384
$ws.Range("D23").Value = 1.033345586974014
$ws.Range("E23").Value = 1.024698424228091
$ws.Range("F23").Value = 1.041521254244513
$ws.Range("I23").Value = 1.033073708509836
$ws.Range("J23").Value = 1.030859923977745
$ws.Range("K23").Value = 1.036996503796713
$ws.Range("L23").Value = 1.028382599775118
$ws.Range("M23").Value = 1.045141288443897

$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.025449395317832
$ws.Range("D24").Value = 1.034368189517167
$ws.Range("E24").Value = 1.025864382318433
$ws.Range("F24").Value = 1.042748130136607
$ws.Range("I24").Value = 1.033298613669594
$ws.Range("J24").Value = 1.031710756142654
$ws.Range("K24").Value = 1.037743955312948
$ws.Range("L24").Value = 1.029270198966426
$ws.Range("M24").Value = 1.046094810564027

$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.027052350972482
$ws.Range("D25").Value = 1.03555673289764
$ws.Range("E25").Value = 1.027221851787698
$ws.Range("F25").Value = 1.044175260987338
$ws.Range("I25").Value = 1.03355422507682
$ws.Range("J25").Value = 1.032698339954248
$ws.Range("K25").Value = 1.038610026896909
$ws.Range("L25").Value = 1.030301586055917
$ws.Range("M25").Value = 1.047201699757722
